# Refresh cached Universalis market-price columns (H:N) on each job sheet.
# Values come from the scheduled market-data pull; rows/cells not listed are unchanged.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38: Hi-Potion of Strength
$ws.Cells.Item(38, 8).Value = 2866
$ws.Cells.Item(38, 9).Value = 2239.2
$ws.Cells.Item(38, 10).Value = 6000
$ws.Cells.Item(38, 11).Value = 6717.599999999999
$ws.Cells.Item(38, 12).Value = 18000
$ws.Cells.Item(38, 13).Value = -6345.599999999999
$ws.Cells.Item(38, 14).Value = -18744

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Bronze Ingot
$ws.Cells.Item(2, 8).Value = 3521.2354
$ws.Cells.Item(2, 9).Value = 2450.077
$ws.Cells.Item(2, 11).Value = 2450.077
$ws.Cells.Item(2, 13).Value = -2337.077

# Row 45: Mythril Ingot
$ws.Cells.Item(45, 8).Value = 1984
$ws.Cells.Item(45, 9).Value = 1942.5
$ws.Cells.Item(45, 10).Value = 2011.6666
$ws.Cells.Item(45, 11).Value = 1942.5
$ws.Cells.Item(45, 12).Value = 2011.6666
$ws.Cells.Item(45, 13).Value = -1565.5
$ws.Cells.Item(45, 14).Value = -2765.6666

# Row 63: Mythrite Rivets
$ws.Cells.Item(63, 8).Value = 2016.2667
$ws.Cells.Item(63, 9).Value = 2104.2
$ws.Cells.Item(63, 10).Value = 1972.3
$ws.Cells.Item(63, 11).Value = 2104.2
$ws.Cells.Item(63, 12).Value = 1972.3
$ws.Cells.Item(63, 13).Value = -1418.2
$ws.Cells.Item(63, 14).Value = -3344.3

# Row 66: Mythrite Rivets
$ws.Cells.Item(66, 8).Value = 2016.2667
$ws.Cells.Item(66, 9).Value = 2104.2
$ws.Cells.Item(66, 10).Value = 1972.3
$ws.Cells.Item(66, 11).Value = 10521
$ws.Cells.Item(66, 12).Value = 9861.5
$ws.Cells.Item(66, 13).Value = -7089
$ws.Cells.Item(66, 14).Value = -16725.5

# Row 74: Titanium Nugget
$ws.Cells.Item(74, 8).Value = 5145113.5
$ws.Cells.Item(74, 9).Value = 3087615.8
$ws.Cells.Item(74, 11).Value = 3087615.8
$ws.Cells.Item(74, 13).Value = -3086741.8

# Row 77: Titanium Nugget
$ws.Cells.Item(77, 8).Value = 5145113.5
$ws.Cells.Item(77, 9).Value = 3087615.8
$ws.Cells.Item(77, 11).Value = 15438079
$ws.Cells.Item(77, 13).Value = -15433711

# Row 102: Tama-hagane Ingot
$ws.Cells.Item(102, 8).Value = 2050
$ws.Cells.Item(102, 9).Value = 1862.5
$ws.Cells.Item(102, 11).Value = 1862.5
$ws.Cells.Item(102, 13).Value = -240.5

# Row 116: Titanbronze Ingot
$ws.Cells.Item(116, 8).Value = 3521.2354
$ws.Cells.Item(116, 9).Value = 2450.077
$ws.Cells.Item(116, 11).Value = 2450.077
$ws.Cells.Item(116, 13).Value = -156.0770000000002

# Row 122: High Durium Nugget
$ws.Cells.Item(122, 8).Value = 1658.2307
$ws.Cells.Item(122, 9).Value = 1721.875
$ws.Cells.Item(122, 10).Value = 1556.4
$ws.Cells.Item(122, 11).Value = 5165.625
$ws.Cells.Item(122, 12).Value = 4669.200000000001
$ws.Cells.Item(122, 13).Value = -2715.625
$ws.Cells.Item(122, 14).Value = -9569.200000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Bronze Ingot
$ws.Cells.Item(3, 8).Value = 3521.2354
$ws.Cells.Item(3, 9).Value = 2450.077
$ws.Cells.Item(3, 11).Value = 2450.077
$ws.Cells.Item(3, 13).Value = -2336.077

# Row 99: Oroshigane Ingot
$ws.Cells.Item(99, 8).Value = 4352.8184
$ws.Cells.Item(99, 9).Value = 4215
$ws.Cells.Item(99, 11).Value = 4215
$ws.Cells.Item(99, 13).Value = -2717

# Row 105: Molybdenum Ingot
$ws.Cells.Item(105, 8).Value = 3679.6667
$ws.Cells.Item(105, 9).Value = 3980.4443
$ws.Cells.Item(105, 11).Value = 3980.4443
$ws.Cells.Item(105, 13).Value = -2233.4443

# Row 107: Deepgold Nugget
$ws.Cells.Item(107, 8).Value = 3546.8
$ws.Cells.Item(107, 10).Value = 4998.5
$ws.Cells.Item(107, 12).Value = 4998.5
$ws.Cells.Item(107, 14).Value = -8838.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Ash Lumber
$ws.Cells.Item(16, 8).Value = 999.5
$ws.Cells.Item(16, 9).Value = 999.5
$ws.Cells.Item(16, 11).Value = 999.5
$ws.Cells.Item(16, 13).Value = -712.5

# Row 31: Walnut Lumber
$ws.Cells.Item(31, 8).Value = 2247.375
$ws.Cells.Item(31, 9).Value = 1403.091
$ws.Cells.Item(31, 10).Value = 4104.8
$ws.Cells.Item(31, 11).Value = 1403.091
$ws.Cells.Item(31, 12).Value = 4104.8
$ws.Cells.Item(31, 13).Value = -1108.091
$ws.Cells.Item(31, 14).Value = -4694.8

# Row 34: Walnut Lumber
$ws.Cells.Item(34, 8).Value = 2247.375
$ws.Cells.Item(34, 9).Value = 1403.091
$ws.Cells.Item(34, 10).Value = 4104.8
$ws.Cells.Item(34, 11).Value = 1403.091
$ws.Cells.Item(34, 12).Value = 4104.8
$ws.Cells.Item(34, 13).Value = -1201.091
$ws.Cells.Item(34, 14).Value = -4508.8

# Row 58: Mahogany Lumber
$ws.Cells.Item(58, 8).Value = 3049.6667
$ws.Cells.Item(58, 9).Value = 2699.75
$ws.Cells.Item(58, 11).Value = 2699.75
$ws.Cells.Item(58, 13).Value = -2496.75

# Row 88: Adamantite Spear
$ws.Cells.Item(88, 8).Value = 32068
$ws.Cells.Item(88, 10).Value = 32068
$ws.Cells.Item(88, 12).Value = 32068
$ws.Cells.Item(88, 14).Value = -32880

# Row 91: Adamantite Spear
$ws.Cells.Item(91, 8).Value = 32068
$ws.Cells.Item(91, 10).Value = 32068
$ws.Cells.Item(91, 12).Value = 32068
$ws.Cells.Item(91, 14).Value = -34876

# Row 105: Zelkova Lumber
$ws.Cells.Item(105, 8).Value = 2725.923
$ws.Cells.Item(105, 9).Value = 2111.6365
$ws.Cells.Item(105, 11).Value = 2111.6365
$ws.Cells.Item(105, 13).Value = -364.6365000000001

# Row 113: White Ash Lumber
$ws.Cells.Item(113, 8).Value = 999.5
$ws.Cells.Item(113, 9).Value = 999.5
$ws.Cells.Item(113, 11).Value = 999.5
$ws.Cells.Item(113, 13).Value = 1170.5

# Row 134: Ceiba Lumber
$ws.Cells.Item(134, 8).Value = 7696138
$ws.Cells.Item(134, 9).Value = 4209.1816
$ws.Cells.Item(134, 11).Value = 12627.5448
$ws.Cells.Item(134, 13).Value = -10092.5448

# Row 136: Dark Mahogany Lumber
$ws.Cells.Item(136, 8).Value = 3049.6667
$ws.Cells.Item(136, 9).Value = 2699.75
$ws.Cells.Item(136, 11).Value = 8099.25
$ws.Cells.Item(136, 13).Value = -5549.25

$ws = $wb.Worksheets.Item("CUL")
# Row 8: Sweet Cream
$ws.Cells.Item(8, 8).Value = 631.4545000000001
$ws.Cells.Item(8, 9).Value = 631.4545000000001
$ws.Cells.Item(8, 11).Value = 1894.3635
$ws.Cells.Item(8, 13).Value = -1755.3635

# Row 14: Kukuru Powder
$ws.Cells.Item(14, 8).Value = 120
$ws.Cells.Item(14, 9).Value = 120
$ws.Cells.Item(14, 11).Value = 360
$ws.Cells.Item(14, 13).Value = -187

# Row 23: Lavender Oil
$ws.Cells.Item(23, 8).Value = 1948
$ws.Cells.Item(23, 9).Value = 2996
$ws.Cells.Item(23, 11).Value = 8988
$ws.Cells.Item(23, 13).Value = -8753

# Row 76: Dhalmel Fricassee
$ws.Cells.Item(76, 8).Value = 9666.666999999999
$ws.Cells.Item(76, 9).Value = 4333.3335
$ws.Cells.Item(76, 10).Value = 15000
$ws.Cells.Item(76, 11).Value = 13000.0005
$ws.Cells.Item(76, 12).Value = 45000
$ws.Cells.Item(76, 13).Value = -12617.0005
$ws.Cells.Item(76, 14).Value = -45766

# Row 79: Dhalmel Fricassee
$ws.Cells.Item(79, 8).Value = 9666.666999999999
$ws.Cells.Item(79, 9).Value = 4333.3335
$ws.Cells.Item(79, 10).Value = 15000
$ws.Cells.Item(79, 11).Value = 13000.0005
$ws.Cells.Item(79, 12).Value = 45000
$ws.Cells.Item(79, 13).Value = -11674.0005
$ws.Cells.Item(79, 14).Value = -47652

# Row 81: Frozen Spirits
$ws.Cells.Item(81, 8).Value = 4642.8887
$ws.Cells.Item(81, 9).Value = 2723.25
$ws.Cells.Item(81, 11).Value = 8169.75
$ws.Cells.Item(81, 13).Value = -7046.75

# Row 84: Frozen Spirits
$ws.Cells.Item(84, 8).Value = 4642.8887
$ws.Cells.Item(84, 9).Value = 2723.25
$ws.Cells.Item(84, 11).Value = 24509.25
$ws.Cells.Item(84, 13).Value = -18893.25

$ws = $wb.Worksheets.Item("GSM")
# Row 97: Koppranickel Ingot
$ws.Cells.Item(97, 8).Value = 1750
$ws.Cells.Item(97, 9).Value = 1750
$ws.Cells.Item(97, 11).Value = 1750
$ws.Cells.Item(97, 13).Value = -1254

# Row 113: Manasilver Nugget
$ws.Cells.Item(113, 8).Value = 2995.4
$ws.Cells.Item(113, 10).Value = 3000
$ws.Cells.Item(113, 12).Value = 3000
$ws.Cells.Item(113, 14).Value = -7340

# Row 132: Lar Ingot
$ws.Cells.Item(132, 8).Value = 2913.8572
$ws.Cells.Item(132, 9).Value = 2899.5
$ws.Cells.Item(132, 11).Value = 8698.5
$ws.Cells.Item(132, 13).Value = -6168.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Boar Leather
$ws.Cells.Item(46, 8).Value = 5428.143
$ws.Cells.Item(46, 10).Value = 5199.6
$ws.Cells.Item(46, 12).Value = 5199.6
$ws.Cells.Item(46, 14).Value = -5575.6

# Row 122: Gaja Leather
$ws.Cells.Item(122, 8).Value = 3419.8
$ws.Cells.Item(122, 10).Value = 3445.0588
$ws.Cells.Item(122, 12).Value = 10335.1764
$ws.Cells.Item(122, 14).Value = -15235.1764

$ws = $wb.Worksheets.Item("WVR")
# Row 45: Linen Trousers
$ws.Cells.Item(45, 8).Value = 18028.875
$ws.Cells.Item(45, 10).Value = 23437.25
$ws.Cells.Item(45, 12).Value = 23437.25
$ws.Cells.Item(45, 14).Value = -24419.25

# Row 113: Pixie Floss
$ws.Cells.Item(113, 8).Value = 3318.7856
$ws.Cells.Item(113, 9).Value = 420.7143
$ws.Cells.Item(113, 11).Value = 1262.1429
$ws.Cells.Item(113, 13).Value = 907.8571000000002

# Row 122: Dark Hempen Cloth
$ws.Cells.Item(122, 8).Value = 4656.533
$ws.Cells.Item(122, 9).Value = 5495.778
$ws.Cells.Item(122, 10).Value = 3397.6667
$ws.Cells.Item(122, 11).Value = 16487.334
$ws.Cells.Item(122, 12).Value = 10193.0001
$ws.Cells.Item(122, 13).Value = -14037.334
$ws.Cells.Item(122, 14).Value = -15093.0001

# Row 136: Sarcenet Cloth
$ws.Cells.Item(136, 8).Value = 946.25
$ws.Cells.Item(136, 9).Value = 948.3333
$ws.Cells.Item(136, 10).Value = 940
$ws.Cells.Item(136, 11).Value = 2844.9999
$ws.Cells.Item(136, 12).Value = 2820
$ws.Cells.Item(136, 13).Value = -294.9998999999998
$ws.Cells.Item(136, 14).Value = -7920

Write-Output "Updated market-price columns on ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
